$wb = $excel.ActiveWorkbook

# Sheet 1: "Branch lengths" - update three branch_length outlier values
$wsBranch = $wb.Worksheets.Item("Branch lengths")
$wsBranch.Range("B97").Value = 1
$wsBranch.Range("B106").Value = 1
$wsBranch.Range("B138").Value = 1

# Sheet 2: "Branch lengths stats" - update recomputed summary statistics
$wsStats = $wb.Worksheets.Item("Branch lengths stats")
$wsStats.Range("B3").Value = 2.537974683544304
$wsStats.Range("B4").Value = 3.364551201745027
$wsStats.Range("B8").Value = 2
